$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RunManager")
$ws2 = $wb.Worksheets.Item("TestData")

# --- RunManager sheet: collapse rows 2-4 into a single "masterModuleTest" row ---
$ws1.Range("A2").Value2 = "masterModuleTest"
$ws1.Range("B2").Value2 = "Yes"
$ws1.Range("C2").Value2 = "'1"
$ws1.Range("D2").Value2 = "'1"

$ws1.Rows(4).Delete()
$ws1.Rows(3).Delete()

# --- View/selection state ---
$ws2.Activate() | Out-Null
$ws2.Range("F1:I1048576").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("E1:H1048576").Select() | Out-Null
